# Auto-generated script applying scheduled market-data refresh to Coeurl_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 4348187.5
$ws.Range("I33").Value = 4545805
$ws.Range("K33").Value = 4545805
$ws.Range("M33").Value = -4545576

$ws.Range("H100").Value = 2546.0833
$ws.Range("I100").Value = 1868.4546
$ws.Range("J100").Value = 10000
$ws.Range("K100").Value = 1868.4546
$ws.Range("L100").Value = 10000
$ws.Range("M100").Value = -1327.4546
$ws.Range("N100").Value = -11082

$ws.Range("H107").Value = 493.82352
$ws.Range("I107").Value = 546.5714
$ws.Range("J107").Value = 247.66667
$ws.Range("K107").Value = 546.5714
$ws.Range("L107").Value = 247.66667
$ws.Range("M107").Value = 1373.4286
$ws.Range("N107").Value = -4087.66667

$ws.Range("H112").Value = 31019.828
$ws.Range("I112").Value = 1990.5883
$ws.Range("J112").Value = 58436.332
$ws.Range("K112").Value = 5971.7649
$ws.Range("L112").Value = 175308.996
$ws.Range("M112").Value = -4863.7649
$ws.Range("N112").Value = -177524.996

$ws.Range("H135").Value = 1294.2727
$ws.Range("I135").Value = 1223.7
$ws.Range("K135").Value = 11013.3
$ws.Range("M135").Value = -8478.300000000001

$ws.Range("H138").Value = 8336418.5
$ws.Range("J138").Value = 11631669
$ws.Range("L138").Value = 34895007
$ws.Range("N138").Value = -34905287

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5550.442
$ws.Range("I74").Value = 1223.3334
$ws.Range("J74").Value = 19829.9
$ws.Range("K74").Value = 1223.3334
$ws.Range("L74").Value = 19829.9
$ws.Range("M74").Value = -349.3334
$ws.Range("N74").Value = -21577.9

$ws.Range("H77").Value = 5550.442
$ws.Range("I77").Value = 1223.3334
$ws.Range("J77").Value = 19829.9
$ws.Range("K77").Value = 6116.666999999999
$ws.Range("L77").Value = 99149.5
$ws.Range("M77").Value = -1748.666999999999
$ws.Range("N77").Value = -107885.5

$ws.Range("H88").Value = 3060
$ws.Range("J88").Value = 3060
$ws.Range("L88").Value = 3060
$ws.Range("N88").Value = -3872

$ws.Range("H91").Value = 3060
$ws.Range("J91").Value = 3060
$ws.Range("L91").Value = 3060
$ws.Range("N91").Value = -5868

$ws.Range("H106").Value = 250000
$ws.Range("J106").Value = 250000
$ws.Range("L106").Value = 250000
$ws.Range("N106").Value = -252524

$ws.Range("H110").Value = 6168
$ws.Range("I110").Value = 6709.35
$ws.Range("K110").Value = 6709.35
$ws.Range("M110").Value = -4664.35

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1293.6207
$ws.Range("I94").Value = 405.5238
$ws.Range("J94").Value = 3624.875
$ws.Range("K94").Value = 405.5238
$ws.Range("L94").Value = 3624.875
$ws.Range("M94").Value = 45.47620000000001
$ws.Range("N94").Value = -4526.875

$ws.Range("H99").Value = 6771.875
$ws.Range("I99").Value = 2362.5
$ws.Range("K99").Value = 2362.5
$ws.Range("M99").Value = -864.5

$ws.Range("H105").Value = 1795.8889
$ws.Range("I105").Value = 1593.8334
$ws.Range("K105").Value = 1593.8334
$ws.Range("M105").Value = 153.1666

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws.Range("H134").Value = 1750.1912
$ws.Range("I134").Value = 1692.5077
$ws.Range("K134").Value = 5077.5231
$ws.Range("M134").Value = -2542.5231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 45855.434
$ws.Range("J31").Value = 6357
$ws.Range("L31").Value = 6357
$ws.Range("N31").Value = -6947

$ws.Range("H34").Value = 45855.434
$ws.Range("J34").Value = 6357
$ws.Range("L34").Value = 6357
$ws.Range("N34").Value = -6761

$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("N56").ClearContents()

$ws.Range("H58").Value = 2372.9375
$ws.Range("I58").Value = 2197.75
$ws.Range("K58").Value = 2197.75
$ws.Range("M58").Value = -1994.75

$ws.Range("H107").Value = 1005.0769
$ws.Range("I107").Value = 558.875
$ws.Range("K107").Value = 558.875
$ws.Range("M107").Value = 1361.125

$ws.Range("H120").Value = 49500
$ws.Range("J120").Value = 49500
$ws.Range("L120").Value = 49500
$ws.Range("N120").Value = -56758

$ws.Range("H122").Value = 1872.5
$ws.Range("I122").Value = 1496.6
$ws.Range("K122").Value = 4489.799999999999
$ws.Range("M122").Value = -2039.799999999999

$ws.Range("H132").Value = 3713.6155
$ws.Range("I132").Value = 3713.6155
$ws.Range("K132").Value = 11140.8465
$ws.Range("M132").Value = -8610.8465

$ws.Range("H134").Value = 33867.3
$ws.Range("I134").Value = 9209.1875
$ws.Range("J134").Value = 132499.75
$ws.Range("K134").Value = 27627.5625
$ws.Range("L134").Value = 397499.25
$ws.Range("M134").Value = -25092.5625
$ws.Range("N134").Value = -402569.25

$ws.Range("H136").Value = 2372.9375
$ws.Range("I136").Value = 2197.75
$ws.Range("K136").Value = 6593.25
$ws.Range("M136").Value = -4043.25
$ws.Range("N136").Value = -13630.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 106256.31
$ws.Range("I14").Value = 106256.31
$ws.Range("K14").Value = 318768.93
$ws.Range("M14").Value = -318595.93

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 64710.25
$ws.Range("I113").Value = 73697.42999999999
$ws.Range("K113").Value = 73697.42999999999
$ws.Range("M113").Value = -71527.42999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6141.5
$ws.Range("I40").Value = 5281.1816
$ws.Range("K40").Value = 5281.1816
$ws.Range("M40").Value = -5145.1816

$ws.Range("H55").Value = 89.94444
$ws.Range("I55").Value = 101.71429
$ws.Range("K55").Value = 101.71429
$ws.Range("M55").Value = 71.28570999999999

$ws.Range("H136").Value = 4675.905
$ws.Range("I136").Value = 4412.933
$ws.Range("K136").Value = 13238.799
$ws.Range("M136").Value = -10688.799

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 997.5
$ws.Range("I96").Value = 997.5
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 997.5
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 375.5
$ws.Range("N96").ClearContents()

$ws.Range("H113").Value = 831.2632
$ws.Range("I113").Value = 831.2632
$ws.Range("K113").Value = 2493.7896
$ws.Range("M113").Value = -323.7896000000001

$ws.Range("H132").Value = 4349.6665
$ws.Range("I132").Value = 4024.5
$ws.Range("K132").Value = 12073.5
$ws.Range("M132").Value = -9543.5

$ws.Range("H136").Value = 2478.889
$ws.Range("I136").Value = 2351.25
$ws.Range("J136").Value = 2843.5715
$ws.Range("K136").Value = 7053.75
$ws.Range("L136").Value = 8530.7145
$ws.Range("M136").Value = -4503.75
$ws.Range("N136").Value = -13630.7145
